$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hotel_info")

# English_Reviews_num, Local_Rank, Total_Reviews_num were blank; fill them in.
# Leading apostrophe forces these numeric-looking values to be stored as
# text (matching the source data's string type), then ClearFormats drops
# the transient "quote prefix" style so the cells keep the default style.
$ws.Range("G2").Value = "'466"
$ws.Range("H2").Value = "'35"
$ws.Range("I2").Value = "'489"
$ws.Range("G2:I2").ClearFormats()
